$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.082
$ws.Range("E3").Value = 12.914
$ws.Range("E5").Value = 13.165
$ws.Range("D9").Value = -8.365
$ws.Range("E11").Value = 13.012
$ws.Range("E12").Value = 13
$ws.Range("D13").Value = -7.662999999999999
$ws.Range("D16").Value = -8.012
$ws.Range("D18").Value = -8.176
$ws.Range("D20").Value = -8.083
$ws.Range("E21").Value = 13.358
